$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '64.913.66'
$ws.Range('E2').Value = '  -1.11%  '

# Row 3
$ws.Range('D3').Value = '3.435.78'
$ws.Range('E3').Value = '  -1.63%  '

# Row 4
$ws.Range('E4').Value = '  -0.01%  '

# Row 5
$ws.Range('D5').Value = "'575.17"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.37%  '

# Row 6
$ws.Range('D6').Value = "'159.63"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.04%  '

# Row 7
$ws.Range('E7').Value = '  +0.01%  '

# Row 8
$ws.Range('B8').Value = 'XRP'
$ws.Range('C8').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D8').Value = "'0.587"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.24%  '

# Row 9
$ws.Range('B9').Value = 'LidoStakedEther'
$ws.Range('C9').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D9').Value = '3.433.85'
$ws.Range('E9').Value = '  -1.83%  '

# Row 10
$ws.Range('D10').Value = "'7.26"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.81%  '

# Row 11
$ws.Range('E11').Value = '  -2.53%  '

# Row 12
$ws.Range('D12').Value = "'0.448"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.13%  '

# Row 13
$ws.Range('D13').Value = '4.027.42'
$ws.Range('E13').Value = '  -1.64%  '

# Row 14
$ws.Range('E14').Value = '  -0.35%  '

# Row 15
$ws.Range('E15').Value = '  -3.50%  '

# Row 16
$ws.Range('D16').Value = "'27.86"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.01%  '

# Row 17
$ws.Range('D17').Value = '64.888.72'
$ws.Range('E17').Value = '  -1.13%  '

# Row 18
$ws.Range('D18').Value = '3.426.57'
$ws.Range('E18').Value = '  -1.98%  '

# Row 19
$ws.Range('D19').Value = "'6.39"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.42%  '

# Row 20
$ws.Range('D20').Value = "'13.93"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.85%  '

# Row 21
$ws.Range('D21').Value = "'382.28"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.73%  '

# Row 22
$ws.Range('D22').Value = "'8.00"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.34%  '

# Row 23
$ws.Range('D23').Value = "'0.550"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.93%  '

# Row 24
$ws.Range('E24').Value = '  +0.37%  '

# Row 25
$ws.Range('D25').Value = "'72.09"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.75%  '

# Row 26
$ws.Range('E26').Value = '  -4.65%  '

# Row 27
$ws.Range('D27').Value = "'9.95"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.11%  '

# Row 28
$ws.Range('E28').Value = '  -0.92%  '

# Row 29
$ws.Range('E29').Value = '  +0.24%  '

# Row 30
$ws.Range('D30').Value = "'1.48"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.77%  '

# Row 31
$ws.Range('D31').Value = "'6.15"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.34%  '

# Row 32
$ws.Range('E32').Value = '  -2.58%  '

# Row 33
$ws.Range('D33').Value = "'23.31"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.82%  '

# Row 34
$ws.Range('D34').Value = "'7.07"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.05%  '

# Row 35
$ws.Range('E35').Value = '  -0.17%  '

# Row 36
$ws.Range('D36').Value = "'160.91"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.08%  '

# Row 37
$ws.Range('D37').Value = "'1.91"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.16%  '

# Row 38
$ws.Range('D38').Value = '2.908.83'
$ws.Range('E38').Value = '  -4.86%  '

# Row 39
$ws.Range('D39').Value = "'0.0752"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.14%  '

# Row 40
$ws.Range('E40').Value = '  +4.06%  '

# Row 41
$ws.Range('D41').Value = "'26.44"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.27%  '

# Row 42
$ws.Range('D42').Value = "'4.59"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.09%  '

# Row 43
$ws.Range('D43').Value = "'43.29"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.30%  '

# Row 44
$ws.Range('D44').Value = "'0.0318"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.24%  '

# Row 45
$ws.Range('E45').Value = '  -0.73%  '

# Row 46
$ws.Range('D46').Value = "'25.99"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.04%  '

# Row 47
$ws.Range('D47').Value = "'2.28"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.14%  '

# Row 48
$ws.Range('D48').Value = "'317.59"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.23%  '

# Row 49
$ws.Range('E49').Value = '  -3.25%  '

# Row 50
$ws.Range('D50').Value = "'6.54"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.17%  '

# Row 51
$ws.Range('E51').Value = '  -2.81%  '
